$wb = $excel.ActiveWorkbook
$ws5 = $wb.Worksheets.Item("GENERADORES")

$ws5.Range("F2:F4").Copy()
$ws5.Range("G2:I4").PasteSpecial(-4122)

$ws5.Range("G1").Value = "RA"
$ws5.Range("H1").Value = "XD"
$ws5.Range("I1").Value = "XQ"
$ws5.Range("G1:I1").Font.ThemeColor = 2
$ws5.Range("G1:I1").Interior.ThemeColor = 1

$ws5.Range("G2").Value = 0.01
$ws5.Range("H2").Value = 2.3
$ws5.Range("I2").Value = 2.2

$ws5.Range("G3").Value = 0.015
$ws5.Range("H3").Value = 2.15
$ws5.Range("I3").Value = 2.23

$ws5.Range("G4").Value = 0.02
$ws5.Range("H4").Value = 2.18
$ws5.Range("I4").Value = 2.25
